$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Commit: "Updated notebook, reran simulation"
#   - sharedStrings: insert "Holden" and "Rizzie Spiral" (two new fiber/method
#     names), rename "Thomas Hex" -> "Matthies Hex"
#   - worksheet: rerun simulation -> new numeric results for every existing
#     row, plus two brand-new rows ("Holden", "Rizzie Spiral") inserted near
#     the top of the table and two more new rows ("Michael-CCHex",
#     "Michael-SNHex") appended at the bottom, growing the table from 28 to
#     30 data rows (dimension A1:T29 -> A1:T31).
# ---------------------------------------------------------------------------

# First, make sure rows 30 and 31 exist with the same formatting as the rest
# of column A (bold, bordered, centered) before we fill in their values --
# copy the format from the last pre-existing data row (29).
$ws.Cells.Item(29,1).Copy() | Out-Null
$ws.Cells.Item(30,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29,1).Copy() | Out-Null
$ws.Cells.Item(31,1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 3: Spiral5
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "Spiral5"
$ws.Cells.Item(3,3).Value = 1.000201973825122
$ws.Cells.Item(3,4).Value = 1.00058283521923
$ws.Cells.Item(3,5).Value = 0.9995516782975582
$ws.Cells.Item(3,6).Value = 1.000169587480294
$ws.Cells.Item(3,7).Value = 1.000169587480294
$ws.Cells.Item(3,8).Value = 0.9992041046004964
$ws.Cells.Item(3,9).Value = 0.9992041046004964
$ws.Cells.Item(3,10).Value = 1.000215298999056
$ws.Cells.Item(3,11).Value = 1.000169587480294
$ws.Cells.Item(3,12).Value = 1.000215298999056
$ws.Cells.Item(3,13).Value = 0.9997097017997763
$ws.Cells.Item(3,14).Value = 0.9997097017997763
$ws.Cells.Item(3,15).Value = 0.9996570272990368
$ws.Cells.Item(3,16).Value = 0.9998629970266154
$ws.Cells.Item(3,17).Value = 0.9998629970266154
$ws.Cells.Item(3,18).Value = 0.999939644640035
$ws.Cells.Item(3,19).Value = 0.999939644640035
$ws.Cells.Item(3,20).Value = 0.9999875797369594

# Row 4: Holden
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Holden"
$ws.Cells.Item(4,3).Value = 1.009766935405545
$ws.Cells.Item(4,4).Value = 1.027944114217353
$ws.Cells.Item(4,5).Value = 0.978341833325804
$ws.Cells.Item(4,6).Value = 1.008377039255956
$ws.Cells.Item(4,7).Value = 1.008377039255956
$ws.Cells.Item(4,8).Value = 0.961447029688834
$ws.Cells.Item(4,9).Value = 0.961447029688834
$ws.Cells.Item(4,10).Value = 1.010338911306067
$ws.Cells.Item(4,11).Value = 1.008377039255956
$ws.Cells.Item(4,12).Value = 1.010338911306067
$ws.Cells.Item(4,13).Value = 0.9858929704974503
$ws.Cells.Item(4,14).Value = 0.9858929704974503
$ws.Cells.Item(4,15).Value = 0.9833759247735682
$ws.Cells.Item(4,16).Value = 0.9933876600836188
$ws.Cells.Item(4,17).Value = 0.9933876600836188
$ws.Cells.Item(4,18).Value = 0.9971350048767031
$ws.Cells.Item(4,19).Value = 0.9971350048767031
$ws.Cells.Item(4,20).Value = 0.9993693105332596

# Row 5: Rizzie Spiral
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Rizzie Spiral"
$ws.Cells.Item(5,3).Value = 1.006749608035177
$ws.Cells.Item(5,4).Value = 1.018984812710395
$ws.Cells.Item(5,5).Value = 0.9850621461142378
$ws.Cells.Item(5,6).Value = 1.00602801634765
$ws.Cells.Item(5,7).Value = 1.00602801634765
$ws.Cells.Item(5,8).Value = 0.9732688322546036
$ws.Cells.Item(5,9).Value = 0.9732688322546036
$ws.Cells.Item(5,10).Value = 1.007046556637662
$ws.Cells.Item(5,11).Value = 1.00602801634765
$ws.Cells.Item(5,12).Value = 1.007046556637662
$ws.Cells.Item(5,13).Value = 0.990157694446133
$ws.Cells.Item(5,14).Value = 0.990157694446133
$ws.Cells.Item(5,15).Value = 0.9884591783355012
$ws.Cells.Item(5,16).Value = 0.9954478017466387
$ws.Cells.Item(5,17).Value = 0.9954478017466387
$ws.Cells.Item(5,18).Value = 0.9980928553968915
$ws.Cells.Item(5,19).Value = 0.9980928553968915
$ws.Cells.Item(5,20).Value = 0.9995233286832876

# Row 6: RotRing OmegaMax-90
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(6,3).Value = 1.003720505617115
$ws.Cells.Item(6,4).Value = 1.009455960600157
$ws.Cells.Item(6,5).Value = 0.9918567738088396
$ws.Cells.Item(6,6).Value = 1.004061134608312
$ws.Cells.Item(6,7).Value = 1.004061134608312
$ws.Cells.Item(6,8).Value = 0.9849918222247689
$ws.Cells.Item(6,9).Value = 0.9849918222247689
$ws.Cells.Item(6,10).Value = 1.003580327990471
$ws.Cells.Item(6,11).Value = 1.004061134608312
$ws.Cells.Item(6,12).Value = 1.003580327990471
$ws.Cells.Item(6,13).Value = 0.9942860751076198
$ws.Cells.Item(6,14).Value = 0.9942860751076198
$ws.Cells.Item(6,15).Value = 0.9934763080080264
$ws.Cells.Item(6,16).Value = 0.9975444282745173
$ws.Cells.Item(6,17).Value = 0.9975444282745173
$ws.Cells.Item(6,18).Value = 0.999173604857966
$ws.Cells.Item(6,19).Value = 0.999173604857966
$ws.Cells.Item(6,20).Value = 0.999611087474944

# Row 7: Equal Angle
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Equal Angle"
$ws.Cells.Item(7,3).Value = 1.011652461851587
$ws.Cells.Item(7,4).Value = 1.033190027038905
$ws.Cells.Item(7,5).Value = 0.9741740676296817
$ws.Cells.Item(7,6).Value = 1.010103122204613
$ws.Cells.Item(7,7).Value = 1.010103122204613
$ws.Cells.Item(7,8).Value = 0.9539639732997099
$ws.Cells.Item(7,9).Value = 0.9539639732997099
$ws.Cells.Item(7,10).Value = 1.012290052536019
$ws.Cells.Item(7,11).Value = 1.010103122204613
$ws.Cells.Item(7,12).Value = 1.012290052536019
$ws.Cells.Item(7,13).Value = 0.9831270129178644
$ws.Cells.Item(7,14).Value = 0.9831270129178644
$ws.Cells.Item(7,15).Value = 0.9801426978218034
$ws.Cells.Item(7,16).Value = 0.9921190493467806
$ws.Cells.Item(7,17).Value = 0.9921190493467806
$ws.Cells.Item(7,18).Value = 0.9966150675612387
$ws.Cells.Item(7,19).Value = 0.9966150675612387
$ws.Cells.Item(7,20).Value = 0.9992289507600861

# Row 8: Tilt Rotate
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "Tilt Rotate"
$ws.Cells.Item(8,3).Value = 1.038912927255321
$ws.Cells.Item(8,4).Value = 1.114210676106502
$ws.Cells.Item(8,5).Value = 0.9134516505879309
$ws.Cells.Item(8,6).Value = 1.031269553004319
$ws.Cells.Item(8,7).Value = 1.031269553004319
$ws.Cells.Item(8,8).Value = 0.847179073998946
$ws.Cells.Item(8,9).Value = 0.847179073998946
$ws.Cells.Item(8,10).Value = 1.042058373649161
$ws.Cells.Item(8,11).Value = 1.031269553004319
$ws.Cells.Item(8,12).Value = 1.042058373649161
$ws.Cells.Item(8,13).Value = 0.9446187238240537
$ws.Cells.Item(8,14).Value = 0.9446187238240537
$ws.Cells.Item(8,15).Value = 0.9342296994120126
$ws.Cells.Item(8,16).Value = 0.9735023335508087
$ws.Cells.Item(8,17).Value = 0.9735023335508087
$ws.Cells.Item(8,18).Value = 0.9879441384141863
$ws.Cells.Item(8,19).Value = 0.9879441384141863
$ws.Cells.Item(8,20).Value = 0.9978470424336967

# Row 9: CLR
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "CLR"
$ws.Cells.Item(9,3).Value = 1.00124068672798
$ws.Cells.Item(9,4).Value = 1.003062326582502
$ws.Cells.Item(9,5).Value = 0.9972926479700712
$ws.Cells.Item(9,6).Value = 1.001420857348934
$ws.Cells.Item(9,7).Value = 1.001420857348934
$ws.Cells.Item(9,8).Value = 0.99497053024916
$ws.Cells.Item(9,9).Value = 0.99497053024916
$ws.Cells.Item(9,10).Value = 1.001166538388408
$ws.Cells.Item(9,11).Value = 1.001420857348934
$ws.Cells.Item(9,12).Value = 1.001166538388408
$ws.Cells.Item(9,13).Value = 0.998068534318784
$ws.Cells.Item(9,14).Value = 0.998068534318784
$ws.Cells.Item(9,15).Value = 0.9978099055358797
$ws.Cells.Item(9,16).Value = 0.999185975328834
$ws.Cells.Item(9,17).Value = 0.999185975328834
$ws.Cells.Item(9,18).Value = 0.9997446958338591
$ws.Cells.Item(9,19).Value = 0.9997446958338591
$ws.Cells.Item(9,20).Value = 0.999858931211176

# Row 10: Rizzie Hex
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "Rizzie Hex"
$ws.Cells.Item(10,3).Value = 1.000098075461255
$ws.Cells.Item(10,4).Value = 1.000266359156377
$ws.Cells.Item(10,5).Value = 0.9997838018786067
$ws.Cells.Item(10,6).Value = 1.000094538121234
$ws.Cells.Item(10,7).Value = 1.000094538121234
$ws.Cells.Item(10,8).Value = 0.999609013479519
$ws.Cells.Item(10,9).Value = 0.999609013479519
$ws.Cells.Item(10,10).Value = 1.00009952926596
$ws.Cells.Item(10,11).Value = 1.000094538121234
$ws.Cells.Item(10,12).Value = 1.00009952926596
$ws.Cells.Item(10,13).Value = 0.9998542713727394
$ws.Cells.Item(10,14).Value = 0.9998542713727394
$ws.Cells.Item(10,15).Value = 0.999830781541362
$ws.Cells.Item(10,16).Value = 0.9999343602889045
$ws.Cells.Item(10,17).Value = 0.9999343602889045
$ws.Cells.Item(10,18).Value = 0.999974404746987
$ws.Cells.Item(10,19).Value = 0.999974404746987
$ws.Cells.Item(10,20).Value = 0.9999918862271585

# Row 11: Matthies Hex
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Matthies Hex"
$ws.Cells.Item(11,3).Value = 1.002194394617832
$ws.Cells.Item(11,4).Value = 1.005402908655259
$ws.Cells.Item(11,5).Value = 0.9952127271496949
$ws.Cells.Item(11,6).Value = 1.002522882961214
$ws.Cells.Item(11,7).Value = 1.002522882961214
$ws.Cells.Item(11,8).Value = 0.9911007653151394
$ws.Cells.Item(11,9).Value = 0.9911007653151394
$ws.Cells.Item(11,10).Value = 1.002059210636925
$ws.Cells.Item(11,11).Value = 1.002522882961214
$ws.Cells.Item(11,12).Value = 1.002059210636925
$ws.Cells.Item(11,13).Value = 0.9965799879760322
$ws.Cells.Item(11,14).Value = 0.9965799879760322
$ws.Cells.Item(11,15).Value = 0.9961242343672531
$ws.Cells.Item(11,16).Value = 0.9985609529710929
$ws.Cells.Item(11,17).Value = 0.9985609529710929
$ws.Cells.Item(11,18).Value = 0.9995514354686232
$ws.Cells.Item(11,19).Value = 0.9995514354686232
$ws.Cells.Item(11,20).Value = 0.9997488148893442

# Row 12: Tilt Rotate_Partial
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(12,3).Value = 1.039483072159759
$ws.Cells.Item(12,4).Value = 1.116056653534117
$ws.Cells.Item(12,5).Value = 0.9121680269759005
$ws.Cells.Item(12,6).Value = 1.031601392428663
$ws.Cells.Item(12,7).Value = 1.031601392428663
$ws.Cells.Item(12,8).Value = 0.8449867517519559
$ws.Cells.Item(12,9).Value = 0.8449867517519559
$ws.Cells.Item(12,10).Value = 1.042726588522507
$ws.Cells.Item(12,11).Value = 1.031601392428663
$ws.Cells.Item(12,12).Value = 1.042726588522507
$ws.Cells.Item(12,13).Value = 0.9438566701372317
$ws.Cells.Item(12,14).Value = 0.9438566701372317
$ws.Cells.Item(12,15).Value = 0.9332937890834546
$ws.Cells.Item(12,16).Value = 0.973104910901042
$ws.Cells.Item(12,17).Value = 0.9731049109010419
$ws.Cells.Item(12,18).Value = 0.9877290312829471
$ws.Cells.Item(12,19).Value = 0.9877290312829471
$ws.Cells.Item(12,20).Value = 0.9978370808954836

# Row 13: RotRing OmegaMax-60
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(13,3).Value = 1.003524584024431
$ws.Cells.Item(13,4).Value = 1.00863416277381
$ws.Cells.Item(13,5).Value = 0.9923147434435825
$ws.Cells.Item(13,6).Value = 1.0040843039854
$ws.Cells.Item(13,7).Value = 1.0040843039854
$ws.Cells.Item(13,8).Value = 0.9856943684391088
$ws.Cells.Item(13,9).Value = 0.9856943684391088
$ws.Cells.Item(13,10).Value = 1.003294244026075
$ws.Cells.Item(13,11).Value = 1.0040843039854
$ws.Cells.Item(13,12).Value = 1.003294244026075
$ws.Cells.Item(13,13).Value = 0.9944943062325918
$ws.Cells.Item(13,14).Value = 0.9944943062325918
$ws.Cells.Item(13,15).Value = 0.993767785302922
$ws.Cells.Item(13,16).Value = 0.9976909721501945
$ws.Cells.Item(13,17).Value = 0.9976909721501945
$ws.Cells.Item(13,18).Value = 0.9992893051089958
$ws.Cells.Item(13,19).Value = 0.9992893051089958
$ws.Cells.Item(13,20).Value = 0.9995910677820677

# Row 14: Equal Angle_Partial
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "Equal Angle_Partial"
$ws.Cells.Item(14,3).Value = 1.011844673789475
$ws.Cells.Item(14,4).Value = 1.034887044957895
$ws.Cells.Item(14,5).Value = 0.9736445948736863
$ws.Cells.Item(14,6).Value = 1.009428410063159
$ws.Cells.Item(14,7).Value = 1.009428410063159
$ws.Cells.Item(14,8).Value = 0.9535161969052628
$ws.Cells.Item(14,9).Value = 0.9535161969052628
$ws.Cells.Item(14,10).Value = 1.012839028221049
$ws.Cells.Item(14,11).Value = 1.009428410063159
$ws.Cells.Item(14,12).Value = 1.012839028221049
$ws.Cells.Item(14,13).Value = 0.9831776125631561
$ws.Cells.Item(14,14).Value = 0.9831776125631561
$ws.Cells.Item(14,15).Value = 0.9799999399999995
$ws.Cells.Item(14,16).Value = 0.9919278783964905
$ws.Cells.Item(14,17).Value = 0.9919278783964905
$ws.Cells.Item(14,18).Value = 0.9963030113131577
$ws.Cells.Item(14,19).Value = 0.9963030113131577
$ws.Cells.Item(14,20).Value = 0.9993599914684212

# Row 15: Rizzie Hex_Partial
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(15,3).Value = 0.992142029787913
$ws.Cells.Item(15,4).Value = 0.9786765933700237
$ws.Cells.Item(15,5).Value = 1.017320726050111
$ws.Cells.Item(15,6).Value = 0.9924119828073844
$ws.Cells.Item(15,7).Value = 0.9924119828073844
$ws.Cells.Item(15,8).Value = 1.031331893897018
$ws.Cells.Item(15,9).Value = 1.031331893897018
$ws.Cells.Item(15,10).Value = 0.992030940415211
$ws.Cells.Item(15,11).Value = 0.9924119828073844
$ws.Cells.Item(15,12).Value = 0.992030940415211
$ws.Cells.Item(15,13).Value = 1.011681417156114
$ws.Cells.Item(15,14).Value = 1.011681417156114
$ws.Cells.Item(15,15).Value = 1.013561186787447
$ws.Cells.Item(15,16).Value = 1.005258272373204
$ws.Cells.Item(15,17).Value = 1.005258272373204
$ws.Cells.Item(15,18).Value = 1.002046699981749
$ws.Cells.Item(15,19).Value = 1.002046699981749
$ws.Cells.Item(15,20).Value = 1.00065236105461

# Row 16: ND Single
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "ND Single"
$ws.Cells.Item(16,3).Value = 1.0681929
$ws.Cells.Item(16,4).Value = 1.200875300000002
$ws.Cells.Item(16,5).Value = 0.8482630399999994
$ws.Cells.Item(16,6).Value = 1.0542662
$ws.Cells.Item(16,7).Value = 1.0542662
$ws.Cells.Item(16,8).Value = 0.7323862499999992
$ws.Cells.Item(16,9).Value = 0.7323862499999992
$ws.Cells.Item(16,10).Value = 1.073924099999999
$ws.Cells.Item(16,11).Value = 1.0542662
$ws.Cells.Item(16,12).Value = 1.073924099999999
$ws.Cells.Item(16,13).Value = 0.9031551749999993
$ws.Cells.Item(16,14).Value = 0.9031551749999993
$ws.Cells.Item(16,15).Value = 0.8848577966666661
$ws.Cells.Item(16,16).Value = 0.9535255166666664
$ws.Cells.Item(16,17).Value = 0.9535255166666664
$ws.Cells.Item(16,18).Value = 0.9787106875
$ws.Cells.Item(16,19).Value = 0.9787106875
$ws.Cells.Item(16,20).Value = 0.9963179650000001

# Row 17: RD Single
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "RD Single"
$ws.Cells.Item(17,3).Value = 1.0262659
$ws.Cells.Item(17,4).Value = 1.0640963
$ws.Cells.Item(17,5).Value = 0.9427502100000001
$ws.Cells.Item(17,6).Value = 1.0306179
$ws.Cells.Item(17,7).Value = 1.0306179
$ws.Cells.Item(17,8).Value = 0.89332461
$ws.Cells.Item(17,9).Value = 0.89332461
$ws.Cells.Item(17,10).Value = 1.0244749
$ws.Cells.Item(17,11).Value = 1.0306179
$ws.Cells.Item(17,12).Value = 1.0244749
$ws.Cells.Item(17,13).Value = 0.958899755
$ws.Cells.Item(17,14).Value = 0.958899755
$ws.Cells.Item(17,15).Value = 0.9535165733333334
$ws.Cells.Item(17,16).Value = 0.9828058033333334
$ws.Cells.Item(17,17).Value = 0.9828058033333332
$ws.Cells.Item(17,18).Value = 0.9947588274999999
$ws.Cells.Item(17,19).Value = 0.9947588274999999
$ws.Cells.Item(17,20).Value = 0.9969216366666668

# Row 18: TD Single
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "TD Single"
$ws.Cells.Item(18,3).Value = 1.0262646
$ws.Cells.Item(18,4).Value = 1.0640928
$ws.Cells.Item(18,5).Value = 0.94275295
$ws.Cells.Item(18,6).Value = 1.0306169
$ws.Cells.Item(18,7).Value = 1.0306169
$ws.Cells.Item(18,8).Value = 0.8933295
$ws.Cells.Item(18,9).Value = 0.8933295
$ws.Cells.Item(18,10).Value = 1.0244736
$ws.Cells.Item(18,11).Value = 1.0306169
$ws.Cells.Item(18,12).Value = 1.0244736
$ws.Cells.Item(18,13).Value = 0.95890155
$ws.Cells.Item(18,14).Value = 0.95890155
$ws.Cells.Item(18,15).Value = 0.9535186833333333
$ws.Cells.Item(18,16).Value = 0.9828066666666667
$ws.Cells.Item(18,17).Value = 0.9828066666666668
$ws.Cells.Item(18,18).Value = 0.9947592250000001
$ws.Cells.Item(18,19).Value = 0.9947592250000001
$ws.Cells.Item(18,20).Value = 0.9969217250000001

# Row 19: Morris Single
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "Morris Single"
$ws.Cells.Item(19,3).Value = 0.9796709
$ws.Cells.Item(19,4).Value = 0.95193969
$ws.Cells.Item(19,5).Value = 1.0441705
$ws.Cells.Item(19,6).Value = 0.97516908
$ws.Cells.Item(19,7).Value = 0.97516908
$ws.Cells.Item(19,8).Value = 1.0829837
$ws.Cells.Item(19,9).Value = 1.0829837
$ws.Cells.Item(19,10).Value = 0.9815235
$ws.Cells.Item(19,11).Value = 0.97516908
$ws.Cells.Item(19,12).Value = 0.9815235
$ws.Cells.Item(19,13).Value = 1.0322536
$ws.Cells.Item(19,14).Value = 1.0322536
$ws.Cells.Item(19,15).Value = 1.0362259
$ws.Cells.Item(19,16).Value = 1.013225426666667
$ws.Cells.Item(19,17).Value = 1.013225426666667
$ws.Cells.Item(19,18).Value = 1.00371134
$ws.Cells.Item(19,19).Value = 1.00371134
$ws.Cells.Item(19,20).Value = 1.002576228333333

# Row 20: Ring Perpendicular to ND
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(20,3).Value = 1.026265480821918
$ws.Cells.Item(20,4).Value = 1.064095184931507
$ws.Cells.Item(20,5).Value = 0.9427511075342468
$ws.Cells.Item(20,6).Value = 1.030617554794521
$ws.Cells.Item(20,7).Value = 1.030617554794521
$ws.Cells.Item(20,8).Value = 0.8933262179452061
$ws.Cells.Item(20,9).Value = 0.8933262179452061
$ws.Cells.Item(20,10).Value = 1.024474475342465
$ws.Cells.Item(20,11).Value = 1.030617554794521
$ws.Cells.Item(20,12).Value = 1.024474475342465
$ws.Cells.Item(20,13).Value = 0.9589003466438355
$ws.Cells.Item(20,14).Value = 0.9589003466438355
$ws.Cells.Item(20,15).Value = 0.9535172669406392
$ws.Cells.Item(20,16).Value = 0.9828060826940638
$ws.Cells.Item(20,17).Value = 0.982806082694064
$ws.Cells.Item(20,18).Value = 0.9947589507191781
$ws.Cells.Item(20,19).Value = 0.9947589507191781
$ws.Cells.Item(20,20).Value = 0.9969216702283106

# Row 21: Ring Perpendicular to RD
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(21,3).Value = 1.011611823157895
$ws.Cells.Item(21,4).Value = 1.033069447368421
$ws.Cells.Item(21,5).Value = 0.9742645636842104
$ws.Cells.Item(21,6).Value = 1.010071437368421
$ws.Cells.Item(21,7).Value = 1.010071437368421
$ws.Cells.Item(21,8).Value = 0.954123197368421
$ws.Cells.Item(21,9).Value = 0.954123197368421
$ws.Cells.Item(21,10).Value = 1.012245740526316
$ws.Cells.Item(21,11).Value = 1.010071437368421
$ws.Cells.Item(21,12).Value = 1.012245740526316
$ws.Cells.Item(21,13).Value = 0.9831844689473684
$ws.Cells.Item(21,14).Value = 0.9831844689473684
$ws.Cells.Item(21,15).Value = 0.9802111671929824
$ws.Cells.Item(21,16).Value = 0.9921467917543859
$ws.Cells.Item(21,17).Value = 0.9921467917543859
$ws.Cells.Item(21,18).Value = 0.9966279531578947
$ws.Cells.Item(21,19).Value = 0.9966279531578947
$ws.Cells.Item(21,20).Value = 0.9992310349122806

# Row 22: Ring Perpendicular to TD
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(22,3).Value = 1.011611417368421
$ws.Cells.Item(22,4).Value = 1.033068419473684
$ws.Cells.Item(22,5).Value = 0.9742654689473684
$ws.Cells.Item(22,6).Value = 1.010070964210527
$ws.Cells.Item(22,7).Value = 1.010070964210527
$ws.Cells.Item(22,8).Value = 0.9541248657894736
$ws.Cells.Item(22,9).Value = 0.9541248657894736
$ws.Cells.Item(22,10).Value = 1.012245349473684
$ws.Cells.Item(22,11).Value = 1.010070964210527
$ws.Cells.Item(22,12).Value = 1.012245349473684
$ws.Cells.Item(22,13).Value = 0.9831851076315788
$ws.Cells.Item(22,14).Value = 0.9831851076315788
$ws.Cells.Item(22,15).Value = 0.980211894736842
$ws.Cells.Item(22,16).Value = 0.9921470598245614
$ws.Cells.Item(22,17).Value = 0.9921470598245614
$ws.Cells.Item(22,18).Value = 0.9966280359210526
$ws.Cells.Item(22,19).Value = 0.9966280359210526
$ws.Cells.Item(22,20).Value = 0.9992310808771929

# Row 23: OffsetFTD
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "OffsetFTD"
$ws.Cells.Item(23,3).Value = 0.9911644975742899
$ws.Cells.Item(23,4).Value = 0.9720494985868959
$ws.Cells.Item(23,5).Value = 1.019833143436857
$ws.Cells.Item(23,6).Value = 0.9943770292357982
$ws.Cells.Item(23,7).Value = 0.9943770292357982
$ws.Cells.Item(23,8).Value = 1.034152202205561
$ws.Cells.Item(23,9).Value = 1.034152202205561
$ws.Cells.Item(23,10).Value = 0.989842461027763
$ws.Cells.Item(23,11).Value = 0.9943770292357982
$ws.Cells.Item(23,12).Value = 0.989842461027763
$ws.Cells.Item(23,13).Value = 1.011997331616662
$ws.Cells.Item(23,14).Value = 1.011997331616662
$ws.Cells.Item(23,15).Value = 1.01460926889006
$ws.Cells.Item(23,16).Value = 1.006123897489707
$ws.Cells.Item(23,17).Value = 1.006123897489707
$ws.Cells.Item(23,18).Value = 1.00318718042623
$ws.Cells.Item(23,19).Value = 1.00318718042623
$ws.Cells.Item(23,20).Value = 1.000236472011194

# Row 24: OffsetATD
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "OffsetATD"
$ws.Cells.Item(24,3).Value = 0.9975594406234848
$ws.Cells.Item(24,4).Value = 0.9953660289563483
$ws.Cells.Item(24,5).Value = 1.005200557745764
$ws.Cells.Item(24,6).Value = 0.9961876684333794
$ws.Cells.Item(24,7).Value = 0.9961876684333794
$ws.Cells.Item(24,8).Value = 1.01027030309366
$ws.Cells.Item(24,9).Value = 1.01027030309366
$ws.Cells.Item(24,10).Value = 0.998123952393022
$ws.Cells.Item(24,11).Value = 0.9961876684333794
$ws.Cells.Item(24,12).Value = 0.998123952393022
$ws.Cells.Item(24,13).Value = 1.004197127743341
$ws.Cells.Item(24,14).Value = 1.004197127743341
$ws.Cells.Item(24,15).Value = 1.004531604410815
$ws.Cells.Item(24,16).Value = 1.001527307973354
$ws.Cells.Item(24,17).Value = 1.001527307973354
$ws.Cells.Item(24,18).Value = 1.00019239808836
$ws.Cells.Item(24,19).Value = 1.00019239808836
$ws.Cells.Item(24,20).Value = 1.00045132520761

# Row 25: OffsetF45
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "OffsetF45"
$ws.Cells.Item(25,3).Value = 0.9911639298574813
$ws.Cells.Item(25,4).Value = 0.9720480601021541
$ws.Cells.Item(25,5).Value = 1.019834369555463
$ws.Cells.Item(25,6).Value = 0.9943764317764056
$ws.Cells.Item(25,7).Value = 0.9943764317764056
$ws.Cells.Item(25,8).Value = 1.034154455444146
$ws.Cells.Item(25,9).Value = 1.034154455444146
$ws.Cells.Item(25,10).Value = 0.9898419146454647
$ws.Cells.Item(25,11).Value = 0.9943764317764056
$ws.Cells.Item(25,12).Value = 0.9898419146454647
$ws.Cells.Item(25,13).Value = 1.011998185044805
$ws.Cells.Item(25,14).Value = 1.011998185044805
$ws.Cells.Item(25,15).Value = 1.014610246548358
$ws.Cells.Item(25,16).Value = 1.006124267288672
$ws.Cells.Item(25,17).Value = 1.006124267288672
$ws.Cells.Item(25,18).Value = 1.003187308410606
$ws.Cells.Item(25,19).Value = 1.003187308410606
$ws.Cells.Item(25,20).Value = 1.000236526896852

# Row 26: OffsetA45
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = "OffsetA45"
$ws.Cells.Item(26,3).Value = 0.9975592918069214
$ws.Cells.Item(26,4).Value = 0.9953656128732151
$ws.Cells.Item(26,5).Value = 1.005200883011216
$ws.Cells.Item(26,6).Value = 0.9961875359280463
$ws.Cells.Item(26,7).Value = 0.9961875359280463
$ws.Cells.Item(26,8).Value = 1.010270888361227
$ws.Cells.Item(26,9).Value = 1.010270888361227
$ws.Cells.Item(26,10).Value = 0.9981238003632598
$ws.Cells.Item(26,11).Value = 0.9961875359280463
$ws.Cells.Item(26,12).Value = 0.9981238003632598
$ws.Cells.Item(26,13).Value = 1.004197344362243
$ws.Cells.Item(26,14).Value = 1.004197344362243
$ws.Cells.Item(26,15).Value = 1.004531857245234
$ws.Cells.Item(26,16).Value = 1.001527408217511
$ws.Cells.Item(26,17).Value = 1.001527408217511
$ws.Cells.Item(26,18).Value = 1.000192440145145
$ws.Cells.Item(26,19).Value = 1.000192440145145
$ws.Cells.Item(26,20).Value = 1.000451335390648

# Row 27: OffsetFRD
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = "OffsetFRD"
$ws.Cells.Item(27,3).Value = 0.9911632994094735
$ws.Cells.Item(27,4).Value = 0.9720464816437654
$ws.Cells.Item(27,5).Value = 1.019835755674433
$ws.Cells.Item(27,6).Value = 0.9943757234837083
$ws.Cells.Item(27,7).Value = 0.9943757234837083
$ws.Cells.Item(27,8).Value = 1.034157020917296
$ws.Cells.Item(27,9).Value = 1.034157020917296
$ws.Cells.Item(27,10).Value = 0.9898413144222012
$ws.Cells.Item(27,11).Value = 0.9943757234837083
$ws.Cells.Item(27,12).Value = 0.9898413144222012
$ws.Cells.Item(27,13).Value = 1.011999167669749
$ws.Cells.Item(27,14).Value = 1.011999167669749
$ws.Cells.Item(27,15).Value = 1.01461136367131
$ws.Cells.Item(27,16).Value = 1.006124686274402
$ws.Cells.Item(27,17).Value = 1.006124686274402
$ws.Cells.Item(27,18).Value = 1.003187445576728
$ws.Cells.Item(27,19).Value = 1.003187445576728
$ws.Cells.Item(27,20).Value = 1.00023659925848

# Row 28: OffsetARD
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = "OffsetARD"
$ws.Cells.Item(28,3).Value = 0.9975591296277325
$ws.Cells.Item(28,4).Value = 0.9953651287184779
$ws.Cells.Item(28,5).Value = 1.005201253793241
$ws.Cells.Item(28,6).Value = 0.996187405494251
$ws.Cells.Item(28,7).Value = 0.996187405494251
$ws.Cells.Item(28,8).Value = 1.010271531119735
$ws.Cells.Item(28,9).Value = 1.010271531119735
$ws.Cells.Item(28,10).Value = 0.9981236204652506
$ws.Cells.Item(28,11).Value = 0.996187405494251
$ws.Cells.Item(28,12).Value = 0.9981236204652506
$ws.Cells.Item(28,13).Value = 1.004197575792493
$ws.Cells.Item(28,14).Value = 1.004197575792493
$ws.Cells.Item(28,15).Value = 1.004532135126076
$ws.Cells.Item(28,16).Value = 1.001527519026412
$ws.Cells.Item(28,17).Value = 1.001527519026412
$ws.Cells.Item(28,18).Value = 1.000192490643372
$ws.Cells.Item(28,19).Value = 1.000192490643372
$ws.Cells.Item(28,20).Value = 1.000451344869781

# Row 29: Gaussian Quadrature
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = "Gaussian Quadrature"
$ws.Cells.Item(29,3).Value = 1.004147366272897
$ws.Cells.Item(29,4).Value = 1.010626025181111
$ws.Cells.Item(29,5).Value = 0.9909148157016981
$ws.Cells.Item(29,6).Value = 1.004464757362321
$ws.Cells.Item(29,7).Value = 1.004464757362321
$ws.Cells.Item(29,8).Value = 0.9832929588389819
$ws.Cells.Item(29,9).Value = 0.9832929588389819
$ws.Cells.Item(29,10).Value = 1.004016754687107
$ws.Cells.Item(29,11).Value = 1.004464757362321
$ws.Cells.Item(29,12).Value = 1.004016754687107
$ws.Cells.Item(29,13).Value = 0.9936548567630443
$ws.Cells.Item(29,14).Value = 0.9936548567630443
$ws.Cells.Item(29,15).Value = 0.9927415097425957
$ws.Cells.Item(29,16).Value = 0.9972581569628032
$ws.Cells.Item(29,17).Value = 0.9972581569628032
$ws.Cells.Item(29,18).Value = 0.9990598070626826
$ws.Cells.Item(29,19).Value = 0.9990598070626826
$ws.Cells.Item(29,20).Value = 0.9995771130073526

# Row 30: Michael-CCHex
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = "Michael-CCHex"
$ws.Cells.Item(30,3).Value = 0.9993521201980604
$ws.Cells.Item(30,4).Value = 0.9969971790118994
$ws.Cells.Item(30,5).Value = 1.001540099107326
$ws.Cells.Item(30,6).Value = 1.000285427124532
$ws.Cells.Item(30,7).Value = 1.000285427124532
$ws.Cells.Item(30,8).Value = 1.002245840557825
$ws.Cells.Item(30,9).Value = 1.002245840557825
$ws.Cells.Item(30,10).Value = 0.99896803925818
$ws.Cells.Item(30,11).Value = 1.000285427124532
$ws.Cells.Item(30,12).Value = 0.99896803925818
$ws.Cells.Item(30,13).Value = 1.000606939908002
$ws.Cells.Item(30,14).Value = 1.000606939908002
$ws.Cells.Item(30,15).Value = 1.000917992974443
$ws.Cells.Item(30,16).Value = 1.000499768980179
$ws.Cells.Item(30,17).Value = 1.000499768980179
$ws.Cells.Item(30,18).Value = 1.000446183516267
$ws.Cells.Item(30,19).Value = 1.000446183516267
$ws.Cells.Item(30,20).Value = 0.9998981175429704

# Row 31: Michael-SNHex
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = "Michael-SNHex"
$ws.Cells.Item(31,3).Value = 0.986785608551549
$ws.Cells.Item(31,4).Value = 0.9630108866588296
$ws.Cells.Item(31,5).Value = 1.029229218937007
$ws.Cells.Item(31,6).Value = 0.9880670334923551
$ws.Cells.Item(31,7).Value = 0.9880670334923551
$ws.Cells.Item(31,8).Value = 1.052382952315143
$ws.Cells.Item(31,9).Value = 1.052382952315143
$ws.Cells.Item(31,10).Value = 0.986258277617614
$ws.Cells.Item(31,11).Value = 0.9880670334923551
$ws.Cells.Item(31,12).Value = 0.986258277617614
$ws.Cells.Item(31,13).Value = 1.019320614966379
$ws.Cells.Item(31,14).Value = 1.019320614966379
$ws.Cells.Item(31,15).Value = 1.022623482956588
$ws.Cells.Item(31,16).Value = 1.008902754475038
$ws.Cells.Item(31,17).Value = 1.008902754475038
$ws.Cells.Item(31,18).Value = 1.003693824229367
$ws.Cells.Item(31,19).Value = 1.003693824229367
$ws.Cells.Item(31,20).Value = 1.00095566292875
